$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the row above (E55:I55) down into the
# newly-populated cells so they pick up the same styles (s="13"/"14").
$ws.Range("E55:I55").Copy()
$ws.Range("E56:I56").PasteSpecial(-4122)
$ws.Range("E57:I57").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 56: "Using ALL to avoid context transition" / ALL, ALLEXCEPT
$ws.Range("E56").Value = 5
$ws.Range("F56").Value = "Using ALL to avoid context transition"
$ws.Range("G56").Value = 1
$ws.Range("H56").Value = "Using ALL to avoid context transition"
$ws.Range("I56").Value = "ALL, ALLEXCEPT"

# Row 57: "Using ISEMPTY" / ISEMPTY
$ws.Range("E57").Value = 6
$ws.Range("F57").Value = "Using ISEMPTY"
$ws.Range("G57").Value = 1
$ws.Range("H57").Value = "Using ISEMPTY"
$ws.Range("I57").Value = "ISEMPTY"

# Update selection/view state to match the committed workbook
$ws.Range("I57").Select()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 2
